$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cols = @("B","C","D","E","F","G","H","I","J","K","L")

# Row 1: forecast date headers (shifted forward one week)
$row1 = @(45685,45692,45699,45706,45713,45720,45727,45734,45741,45748,45755)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "1").Value = $row1[$i] }

# Row 2: _Average
$row = @(61.1,60.4,60,59.9,60.3,61,62.3,63.8,65.7,67.7,69.7)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "2").Value = $row[$i] }

# Row 3: Abim
$row = @(64.09999999999999,65.90000000000001,67.7,69.59999999999999,71.5,73.3,75,76.40000000000001,77.59999999999999,78.5,78.90000000000001)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "3").Value = $row[$i] }

# Row 4: Amudat
$row = @(42.5,41.8,41.5,41.6,42.3,43.6,45.4,47.7,50.4,53.3,56.3)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "4").Value = $row[$i] }

# Row 5: Budi
$row = @(64.2,62.6,61.2,60.2,59.5,59.2,59.4,60,61,62.2,63.6)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "5").Value = $row[$i] }

# Row 6: Dasenech (Kuraz)
$row = @(64.59999999999999,66.7,68.90000000000001,71,73,74.7,76.2,77.3,78.2,78.7,78.90000000000001)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "6").Value = $row[$i] }

# Row 7: Gnangatom
$row = @(60.4,59.9,59.7,59.7,60.2,61.3,62.9,65.09999999999999,67.59999999999999,70.40000000000001,73.2)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "7").Value = $row[$i] }

# Row 8: Kaabong
$row = @(68.59999999999999,67.09999999999999,65.8,65,64.59999999999999,64.7,65.3,66.40000000000001,68,69.90000000000001,71.90000000000001)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "8").Value = $row[$i] }

# Row 9: Kapoeta East
$row = @(60.8,58.7,57.1,56.1,55.8,56.4,57.8,60.1,62.9,66.09999999999999,69.5)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "9").Value = $row[$i] }

# Row 10: Kapoeta North
$row = @(60.7,57.9,55.5,53.8,52.9,52.7,53.4,54.8,56.8,59.3,62)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "10").Value = $row[$i] }

# Row 11: Kapoeta South
$row = @(72.59999999999999,70.40000000000001,68.5,67,66.09999999999999,65.90000000000001,66.3,67.40000000000001,69.09999999999999,71.09999999999999,73.3)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "11").Value = $row[$i] }

# Row 12: Karenga
$row = @(65.3,63.8,62.3,60.9,59.5,58.3,57.3,56.5,55.9,55.6,55.4)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "12").Value = $row[$i] }

# Row 13: Kotido
$row = @(62.3,61.4,60.6,60,59.7,59.7,60,60.5,61.2,62,62.8)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "13").Value = $row[$i] }

# Row 14: Loima
$row = @(58.9,60.2,62.6,65.90000000000001,70.2,75.40000000000001,81.3,87.5,93.90000000000001,100,105.4)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "14").Value = $row[$i] }

# Row 15: Moroto
$row = @(62,59.3,56.8,54.9,53.6,53.2,53.8,55.2,57.4,60.3,63.5)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "15").Value = $row[$i] }

# Row 16: Nakapiripirit
$row = @(55.2,54.4,53.9,53.7,53.7,53.8,54.2,54.8,55.5,56.2,56.9)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "16").Value = $row[$i] }

# Row 17: Napak
$row = @(62.4,61.4,60.5,59.7,59.1,58.7,58.5,58.6,58.8,59.3,59.9)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "17").Value = $row[$i] }

# Row 18: Surma
$row = @(61.5,59.6,58,56.9,56.4,56.4,57.2,58.5,60.4,62.5,64.90000000000001)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "18").Value = $row[$i] }

# Row 19: Turkana
$row = @(56.9,58.8,61.2,63.8,66.59999999999999,69.59999999999999,72.59999999999999,75.59999999999999,78.40000000000001,80.90000000000001,82.90000000000001)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "19").Value = $row[$i] }

# Row 20: Turkana West
$row = @(64.7,65,65.59999999999999,66.5,67.8,69.59999999999999,71.90000000000001,74.7,77.8,81,84.09999999999999)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "20").Value = $row[$i] }

# Row 21: West Pokot
$row = @(53.3,52.8,52.5,52.5,52.8,53.4,54.3,55.4,56.8,58.4,60)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "21").Value = $row[$i] }
